# Add the new "Lords Stronghold" location row (row 57) to the Locations sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 57

$ws.Range("A$row").Value = "Lords Stronghold"
$ws.Range("B$row").Value = "Shadow Plane"
$ws.Range("D$row").Value = "Key to the Stronghold"
$ws.Range("E$row").Value = "The Lords of the Shade Realm live with in this strong hold. Seeking nothing but total control of the past."
$ws.Range("G$row").Value = 1
$ws.Range("J$row").Value = 2112
$ws.Range("K$row").Value = 2112
$ws.Range("L$row").Value = 8
$ws.Range("M$row").Value = "No"
